$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("J1").Value = "TS estimate"
$ws.Range("K1").Value = "TS Actual"

# Data rows
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 60
$ws.Range("J4").Value = 12
$ws.Range("K4").Value = 12
$ws.Range("J21").Value = 4
$ws.Range("K21").Value = 4
$ws.Range("J24").Value = 8
$ws.Range("K24").Value = 0.5
$ws.Range("J25").Value = 8
$ws.Range("K25").Value = 0.75
$ws.Range("J26").Value = 8
$ws.Range("K26").Value = 0.75
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 1
$ws.Range("J28").Value = 16
$ws.Range("K28").Value = 1
$ws.Range("J29").Value = 16
$ws.Range("K29").Value = 0.75
$ws.Range("J30").Value = 8
$ws.Range("K30").Value = 1
$ws.Range("J31").Value = 8
$ws.Range("K31").Value = 0.75
$ws.Range("J32").Value = 8
$ws.Range("K32").Value = 2
$ws.Range("J33").Value = 8
$ws.Range("K33").Value = 2
$ws.Range("J34").Value = 8
$ws.Range("K34").Value = 1
$ws.Range("J35").Value = 8
$ws.Range("K35").Value = 2
$ws.Range("J37").Value = 16
$ws.Range("K37").Value = 2

# Sum row - copy the shared formula across
$ws.Range("I39:K39").FillRight()

# New ratio row
$ws.Range("J43").Formula = "=K39/J39"
$ws.Range("J43").NumberFormat = "0%"
